# Sprint 1 Meeting 8 - populate H/I columns (meeting 8 Q&A) and adjust layout
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values for columns H (3/11/2019: 6pm) and I (3/13/2019: 6pm) ---
$ws.Range("H2").Value = 'Created/updated layouts for user settings, change password login, forgot, and set security questions pages. Used logos that the team agreed upon in layouts.'
$ws.Range("I2").Value = 'Moved project husk into git repo with several screen layouts finished'
$ws.Range("H3").Value = 'I worked on the security issue where the user is signed out after 5 minutes of inactivty and the startup page issue and the create account issue, and I created more Mazes for future levels'
$ws.Range("I3").Value = 'I worked on implementation of the startup page issue and the "create account" issue, and I created more mazes for future levels'
$ws.Range("H4").Value = 'No significant progress'
$ws.Range("I4").Value = 'No significant progress'
$ws.Range("H5").Value = 'Completed animation for the sprite for the main character for the Hack Heist app and began work on the user settings button/screen for the in-game menu'
$ws.Range("I5").Value = 'I have began working on exporting Unity projects to android studio'
$ws.Range("H6").Value = 'Looking into intentions and utilizing Google Firebase to control data'
$ws.Range("I6").Value = 'Looking into intentions and utilizing Google Firebase to control data'
$ws.Range("H7").Value = 'I will work on my assigned issues. Every issue assigned to me will be worked on as well as more level design.'
$ws.Range("I7").Value = 'I will work on my assigned issues. Every issue assigned to me will be worked on as well as more level design.'
$ws.Range("H8").Value = 'Further work on assigned pages'
$ws.Range("I8").Value = 'Further work on assigned pages'
$ws.Range("H9").Value = 'Continue work on the user settings button/screen for the in-game menu. Begin work on the next sprite for the game.'
$ws.Range("I9").Value = 'Continue attempting to export projects from Unity to Android studio'
$ws.Range("H10").Value = 'Not currently'
$ws.Range("I10").Value = 'Not currently'
$ws.Range("H11").Value = 'No, I am finally working with all cylinders pumping'
$ws.Range("I11").Value = 'No, I am working with all cylinders pumping'
$ws.Range("H12").Value = 'London'
$ws.Range("I12").Value = 'London'
$ws.Range("H13").Value = 'Nothing is currently getting in the way of my work.'
$ws.Range("I13").Value = 'Nothing is currently getting in the way of my work.'
$ws.Range("H14").Value = 'Nothing unexpected as of yet'
$ws.Range("I14").Value = 'Nothing unexpected as of yet'
$ws.Range("H15").Value = 'It takes time to get into a good rhythm, but once found, ride it!'
$ws.Range("I15").Value = 'Documentation can really slow down the process of good work, but perhaps has hidden value to be discovered'
$ws.Range("H16").Value = 'Connecting activities/post spring break I need to spend a lot more time on the project'
$ws.Range("I16").Value = 'Post spring break I need to spend a lot more time on the project'
$ws.Range("H17").Value = 'Currently learning how to make and connect multiple activities in android studio'
$ws.Range("I17").Value = 'Currently learning how to modify and build settings on the Unity project'
$ws.Range("H18").Value = 'Not currently'
$ws.Range("I18").Value = 'Not currently'
$ws.Range("H19").Value = 'I thing that our group should be more on the same page of the design of the app I think we should follow the Model|Presenter|View method for software development. However, we are just sorting hacking issues individually and at the end we are supposidly going to mush all of our work into one project. I am concerned about this getting out of control later on. '
$ws.Range("I19").Value = 'I thing that our group should be more on the same page of the design of the app I think we should follow the Model|Presenter|View method for software development. However, we are just sorting hacking issues individually and at the end we are supposidly going to mush all of our work into one project. I am concerned about this getting out of control later on. '
$ws.Range("H20").Value = 'Not yet'
$ws.Range("I20").Value = 'Not yet'
$ws.Range("H21").Value = 'No changes currently need to be made to the project due to my work'
$ws.Range("I21").Value = 'No changes currently need to be made to the project due to my work'

# --- Row 19 (Q: "Do any changes to the project need to be made...") uses the smaller/condensed font + wrap style ---
$ws.Range("H19").Font.Size = 9
$ws.Range("H19").WrapText = $true
$ws.Range("I19").Font.Size = 9
$ws.Range("I19").WrapText = $true

# --- Row heights adjusted to fit the new, longer answers ---
$ws.Rows.Item(3).RowHeight = 108
$ws.Rows.Item(5).RowHeight = 97
$ws.Rows.Item(9).RowHeight = 81.5
$ws.Rows.Item(15).RowHeight = 73.5
$ws.Rows.Item(19).RowHeight = 157

# --- View: scroll position + active selection moved to J19 ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("J19").Select() | Out-Null

